$d = $word.ActiveDocument

$replacements = @(
    @{ old = "522×3="; new = "201×6=" },
    @{ old = "319×9="; new = "937×6=" },
    @{ old = "514×3="; new = "274×7=" },
    @{ old = "967×5="; new = "250×8=" },
    @{ old = "766×2="; new = "947×8=" },
    @{ old = "421×9="; new = "718×8=" },
    @{ old = "541×3="; new = "448×6=" },
    @{ old = "242×6="; new = "625×8=" },
    @{ old = "466×4="; new = "713×2=" },
    @{ old = "105×4="; new = "973×5=" },
    @{ old = "823×2="; new = "113×7=" },
    @{ old = "537×2="; new = "984×6=" },
    @{ old = "782×5="; new = "349×5=" },
    @{ old = "706×7="; new = "535×2=" },
    @{ old = "910×8="; new = "803×5=" },
    @{ old = "632×3="; new = "360×9=" },
    @{ old = "860×7="; new = "589×2=" },
    @{ old = "476×2="; new = "647×8=" },
    @{ old = "939×9="; new = "920×3=" },
    @{ old = "484×9="; new = "633×5=" },
    @{ old = "142×8="; new = "291×9=" },
    @{ old = "806×9="; new = "418×3=" },
    @{ old = "970×5="; new = "624×3=" },
    @{ old = "500×4="; new = "214×3=" },
    @{ old = "624×8="; new = "786×9=" }
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
